# menu bar and  bug fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data corrections
# B2: "fe" -> "48" (keep as text, not a number, so prefix with an apostrophe
# like a real user would to force text entry in a General-formatted cell)
$ws.Range("B2").Value = "'48"
$ws.Range("C2").Value = 33
$ws.Range("D2").Value = "4833@gmail.com"
$ws.Range("E2").Value = "Realtime"
$ws.Range("F2").Value = "['#Volunteering', '#YouthProjects', '#Career', '#Education']"
